$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 192.125
$ws.Range("I2").Value = 192.125
$ws.Range("K2").Value = 192.125
$ws.Range("M2").Value = -79.125

$ws.Range("H18").Value = 6085.7144
$ws.Range("I18").Value = 6085.7144
$ws.Range("K18").Value = 6085.7144
$ws.Range("M18").Value = -5801.7144

$ws.Range("H88").Value = 427263.25
$ws.Range("I88").Value = 1845.4286
$ws.Range("J88").Value = 675423.7
$ws.Range("K88").Value = 1845.4286
$ws.Range("L88").Value = 675423.7
$ws.Range("M88").Value = -1439.4286
$ws.Range("N88").Value = -676235.7

$ws.Range("H91").Value = 427263.25
$ws.Range("I91").Value = 1845.4286
$ws.Range("J91").Value = 675423.7
$ws.Range("K91").Value = 1845.4286
$ws.Range("L91").Value = 675423.7
$ws.Range("M91").Value = -441.4286
$ws.Range("N91").Value = -678231.7

$ws.Range("H92").Value = 1075
$ws.Range("I92").Value = 628.2143
$ws.Range("J92").Value = 2117.5
$ws.Range("K92").Value = 628.2143
$ws.Range("L92").Value = 2117.5
$ws.Range("M92").Value = 619.7857
$ws.Range("N92").Value = -4613.5

$ws.Range("H125").Value = 900

$ws.Range("H138").Value = 3695.375
$ws.Range("J138").Value = 4006.772
$ws.Range("L138").Value = 12020.316
$ws.Range("N138").Value = -22300.316

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 9622749
$ws.Range("I74").Value = 14707004
$ws.Range("K74").Value = 14707004
$ws.Range("M74").Value = -14706130

$ws.Range("H77").Value = 9622749
$ws.Range("I77").Value = 14707004
$ws.Range("K77").Value = 73535020
$ws.Range("M77").Value = -73530652

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2197.375
$ws.Range("I20").Value = 2149
$ws.Range("J20").Value = 2245.75
$ws.Range("K20").Value = 2149
$ws.Range("L20").Value = 2245.75
$ws.Range("M20").Value = -1902
$ws.Range("N20").Value = -2739.75

$ws.Range("H134").Value = 27556.781
$ws.Range("I134").Value = 3287.7693
$ws.Range("K134").Value = 9863.3079
$ws.Range("M134").Value = -7328.3079

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 15404.4
$ws.Range("I16").Value = 15404.4
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 15404.4
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -15117.4
$ws.Range("N16").ClearContents()

$ws.Range("H31").Value = 760340.5600000001
$ws.Range("I31").Value = 12173.4
$ws.Range("K31").Value = 12173.4
$ws.Range("M31").Value = -11878.4

$ws.Range("H34").Value = 760340.5600000001
$ws.Range("I34").Value = 12173.4
$ws.Range("K34").Value = 12173.4
$ws.Range("M34").Value = -11971.4

$ws.Range("H86").Value = 3531.1
$ws.Range("I86").Value = 3488.4
$ws.Range("K86").Value = 3488.4
$ws.Range("M86").Value = -2365.4

$ws.Range("H89").Value = 3531.1
$ws.Range("I89").Value = 3488.4
$ws.Range("K89").Value = 17442
$ws.Range("M89").Value = -11826

$ws.Range("H99").Value = 2995
$ws.Range("I99").Value = 2995
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 2995
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -1497
$ws.Range("N99").ClearContents()

$ws.Range("H113").Value = 15404.4
$ws.Range("I113").Value = 15404.4
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 15404.4
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -13234.4
$ws.Range("N113").ClearContents()

$ws.Range("H126").Value = 2995
$ws.Range("I126").Value = 2995
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 8985
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -6515
$ws.Range("N126").ClearContents()

$ws.Range("H132").Value = 2609.2188
$ws.Range("I132").Value = 2449.8333
$ws.Range("K132").Value = 7349.499899999999
$ws.Range("M132").Value = -4819.499899999999

$ws.Range("H134").Value = 502561.84
$ws.Range("I134").Value = 716324.3
$ws.Range("K134").Value = 2148972.9
$ws.Range("M134").Value = -2146437.9

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H88").Value = 3949.9
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 3949.9
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 11849.7
$ws.Range("N88").Value = -12705.7
$ws.Range("M88").ClearContents()

$ws.Range("H91").Value = 3949.9
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 3949.9
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 11849.7
$ws.Range("N91").Value = -14813.7
$ws.Range("M91").ClearContents()

$ws.Range("H122").Value = 617.5185
$ws.Range("I122").Value = 590.7273
$ws.Range("J122").Value = 635.9375
$ws.Range("K122").Value = 5316.545700000001
$ws.Range("L122").Value = 5723.4375
$ws.Range("M122").Value = -2866.545700000001
$ws.Range("N122").Value = -10623.4375

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2910.2646
$ws.Range("I102").Value = 2438.25
$ws.Range("K102").Value = 2438.25
$ws.Range("M102").Value = -816.25

$ws.Range("H132").Value = 25002934
$ws.Range("I132").Value = 27780726
$ws.Range("J132").Value = 2798.75
$ws.Range("K132").Value = 83342178
$ws.Range("L132").Value = 8396.25
$ws.Range("M132").Value = -83339648
$ws.Range("N132").Value = -13456.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4087.4075
$ws.Range("I40").Value = 3774
$ws.Range("K40").Value = 3774
$ws.Range("M40").Value = -3638

$ws.Range("H55").Value = 83334020
$ws.Range("I55").Value = 142857810
$ws.Range("K55").Value = 142857810
$ws.Range("M55").Value = -142857637

$ws.Range("H93").Value = 71430610
$ws.Range("I93").Value = 111112460
$ws.Range("J93").Value = 3260
$ws.Range("K93").Value = 111112460
$ws.Range("L93").Value = 3260
$ws.Range("M93").Value = -111111212
$ws.Range("N93").Value = -5756

$ws.Range("H122").Value = 6089
$ws.Range("I122").Value = 5435.5625
$ws.Range("J122").Value = 8180
$ws.Range("K122").Value = 16306.6875
$ws.Range("L122").Value = 24540
$ws.Range("M122").Value = -13856.6875
$ws.Range("N122").Value = -29440

$ws.Range("H132").Value = 1004645
$ws.Range("I132").Value = 1433321.4
$ws.Range("K132").Value = 4299964.199999999
$ws.Range("M132").Value = -4297434.199999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 3769.52
$ws.Range("I126").Value = 2842.3635
$ws.Range("J126").Value = 5569.294
$ws.Range("K126").Value = 8527.0905
$ws.Range("L126").Value = 16707.882
$ws.Range("M126").Value = -6057.0905
$ws.Range("N126").Value = -21647.882

$ws.Range("H132").Value = 225911.84
$ws.Range("I132").Value = 2774.3447
$ws.Range("K132").Value = 8323.034100000001
$ws.Range("M132").Value = -5793.034100000001
